$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.078.25"
$ws.Range("E2").Value = "  +0.13%  "

$ws.Range("D3").Value = "1.875.80"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'313.54"
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").Value = "'1.003"
$ws.Range("E6").Value = "  +0.16%  "

$ws.Range("D7").Value = "'0.5071"
$ws.Range("E7").Value = "  -0.12%  "

$ws.Range("D8").Value = "'0.3845"
$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").Value = "'0.08971"
$ws.Range("E9").Value = "  -3.27%  "

$ws.Range("D10").Value = "'1.122"
$ws.Range("E10").Value = "  -0.96%  "

$ws.Range("D11").Value = "'41.58"
$ws.Range("E11").Value = "  -0.61%  "

$ws.Range("D12").Value = "'6.336"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("E13").Value = "  -0.25%  "

$ws.Range("D14").Value = "1.872.39"
$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("D15").Value = "'7.206"
$ws.Range("E15").Value = "  -1.20%  "

$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("D18").Value = "'91.12"
$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("D19").Value = "'0.06599"
$ws.Range("E19").Value = "  +0.30%  "

$ws.Range("D20").Value = "'18.12"
$ws.Range("E20").Value = "  +1.93%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "'6.108"
$ws.Range("E22").Value = "  -1.77%  "

$ws.Range("D23").Value = "28.100.71"
$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("D25").Value = "'2.278"
$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("D26").Value = "2.091.96"
$ws.Range("E26").Value = "  -1.25%  "

$ws.Range("D27").Value = "'2.537"
$ws.Range("E27").Value = "  -2.35%  "

$ws.Range("D28").Value = "'20.74"
$ws.Range("E28").Value = "  -0.88%  "

$ws.Range("D29").Value = "'156.98"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").Value = "'126.60"
$ws.Range("E30").Value = "  -0.49%  "

$ws.Range("E31").Value = "  -1.57%  "

$ws.Range("D32").Value = "'1.060"
$ws.Range("E32").Value = "  -2.23%  "

$ws.Range("E33").Value = "  +0.03%  "

$ws.Range("D34").Value = "'3.603"
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("D35").Value = "'9.626"
$ws.Range("E35").Value = "  +0.45%  "

$ws.Range("D36").Value = "'0.06575"
$ws.Range("E36").Value = "  -1.22%  "

$ws.Range("D37").Value = "'0.02418"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("D38").Value = "'0.2176"
$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").Value = "'1.267"
$ws.Range("E39").Value = "  +1.09%  "

$ws.Range("D40").Value = "'1.206"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("D41").Value = "'0.6393"
$ws.Range("E41").Value = "  +0.69%  "

$ws.Range("D42").Value = "'11.46"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("D43").Value = "'4.909"
$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("D44").Value = "'0.6025"
$ws.Range("E44").Value = "  +0.81%  "

$ws.Range("D45").Value = "'13.16"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("D46").Value = "'3.672"
$ws.Range("E46").Value = "  -0.91%  "

$ws.Range("E47").Value = "  +0.17%  "

$ws.Range("D48").Value = "'1.236"
$ws.Range("E48").Value = "  +5.01%  "

$ws.Range("D49").Value = "'1.993"
$ws.Range("E49").Value = "  -0.74%  "

$ws.Range("D50").Value = "'121.23"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("D51").Value = "'79.69"
$ws.Range("E51").Value = "  +1.91%  "
